# 16.6.1.xlsx update: refresh year columns (2018-2021) with new actual/approved
# figures, add 2021 утв./2021 факт/откл. columns (AK:AM), and move selection to AF4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 3 header relabeling for existing columns (AC, AF, AH, AI, AJ)
# ---------------------------------------------------------------------------
$ws.Range("AC3").Value = "2018 факт "
$ws.Range("AF3").Value = "2019 факт"
$ws.Range("AH3").Value = "2020 утв."
$ws.Range("AI3").Value = "2020 факт"
$ws.Range("AJ3").Value = "откл. от утв., %"

# ---------------------------------------------------------------------------
# 2. Updated data values for the existing 2019/2020 block (columns AE:AJ),
#    rows 5-12 (one data row per budget sector).
# ---------------------------------------------------------------------------
$dataUpdates = @{
    5  = @{ AE = 43737.8;  AF = 43258.3;              AG = 98.9;  AH = 46293.5; AI = 47153.5;              AJ = 101.9 }
    6  = @{ AE = 6265.4;   AF = 4434.6000000000004;    AG = 70.8;  AH = 7935.8;  AI = 3895.8;               AJ = 49.1 }
    7  = @{ AE = 728.5;    AF = 695.7;                 AG = 95.5;  AH = 746.9;   AI = 583.20000000000005;   AJ = 78.099999999999994 }
    8  = @{ AE = 1249;     AF = 1244.7;                AG = 99.7;  AH = 1249;    AI = 1207.5999999999999;   AJ = 96.7 }
    9  = @{ AE = 2582.6;   AF = 2477.5;                AG = 95.9;  AH = 3109;    AI = 3225.2;               AJ = 103.7 }
    10 = @{ AE = 2686.4;   AF = 2829;                  AG = 105.3; AH = 2993.4;  AI = 2624.5;               AJ = 87.7 }
    11 = @{ AE = 23397.4;  AF = 24364.799999999999;    AG = 104.1; AH = 30085.9; AI = 29223.5;              AJ = 97.1 }
    12 = @{ AE = 13137.1;  AF = 10924.7;               AG = 83.2;  AH = 12158.7; AI = 10980.3;              AJ = 90.3 }
}

foreach ($row in $dataUpdates.Keys) {
    $vals = $dataUpdates[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}

# ---------------------------------------------------------------------------
# 3. New columns AK:AM (2021 утв. / 2021 факт / откл. от утв., %).
#    Copy formatting across from the neighbouring AJ column first (so
#    borders/number formats/styles match), then write the actual values.
# ---------------------------------------------------------------------------
$ws.Range("AJ2").Copy()
$ws.Range("AK2:AM2").PasteSpecial(-4122)

$ws.Range("AJ3").Copy()
$ws.Range("AK3:AM3").PasteSpecial(-4122)

$ws.Range("AJ5:AJ12").Copy()
$ws.Range("AK5:AK12").PasteSpecial(-4122)
$ws.Range("AL5:AL12").PasteSpecial(-4122)
$ws.Range("AM5:AM12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 3 headers for the new columns
$ws.Range("AK3").Value = "2021 утв."
$ws.Range("AL3").Value = "2021 факт"
$ws.Range("AM3").Value = "откл. от утв., %"

# New data values for AK:AM, rows 5-12
$newData = @{
    5  = @{ AK = 47483.3;              AL = 52020.5;              AM = 109.6 }
    6  = @{ AK = 8997;                 AL = 6212.4;               AM = 69 }
    7  = @{ AK = 639.20000000000005;   AL = 600.79999999999995;   AM = 94 }
    8  = @{ AK = 1208.0999999999999;   AL = 1332.7;               AM = 110.3 }
    9  = @{ AK = 3131.3;               AL = 4833.7;               AM = 154.4 }
    10 = @{ AK = 2798.4;               AL = 3088;                 AM = 110.3 }
    11 = @{ AK = 30439.7;              AL = 30705.3;              AM = 100.9 }
    12 = @{ AK = 11664.9;              AL = 11939.1;              AM = 102.4 }
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}

# ---------------------------------------------------------------------------
# 4. Move the active selection to AF4, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("AF4").Select()
